# IAM test data workbook - add two new authorize-API test rows (Facebook, LinkedIn)
# and rename the existing TR-ID row's description/fix up column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 1 header: B1 now wraps like the other header cells ---
$ws.Range("B1").WrapText = $true

# --- Row 2: update description text, swap D2/G2 values ---
$ws.Range("B2").Value = "Test the authorize API for Redirection to TR ID login page"
$ws.Range("B2").WrapText = $true
$ws.Range("D2").Value = "/authorize"
$ws.Range("G2").Value = "?provider=thomsonreuters&backurl=%2Fui%2F%23%2Flogin%2FAUTHTOKEN"
$ws.Rows.Item(2).RowHeight = 45

# --- Row 3: new Facebook authorize test case ---
$ws.Range("A3").Value = "S1_TC_T2"
$ws.Range("B3").Value = "Test the authorize API for Redirection to Facebook login page"
$ws.Range("B3").WrapText = $true
$ws.Range("C3").Value = "1PAUTH"
$ws.Range("D3").Value = "/authorize"
$ws.Range("D3").NumberFormat = $ws.Range("D2").NumberFormat
$ws.Range("E3").Value = "GET"
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = "?provider=facebook&backurl=%2Fui%2F%23%2Flogin%2FAUTHTOKEN"
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = "status=200"
$ws.Range("K3").Value = ""
$ws.Rows.Item(3).RowHeight = 30

# --- Row 4: new LinkedIn authorize test case ---
$ws.Range("A4").Value = "S1_TC_T3"
$ws.Range("B4").Value = "Test the authorize API for Redirection to Linked-In login page"
$ws.Range("B4").WrapText = $true
$ws.Range("C4").Value = "1PAUTH"
$ws.Range("D4").Value = "/authorize"
$ws.Range("E4").Value = "GET"
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = "?provider=linkedin&backurl=%2Fui%2F%23%2Flogin%2FAUTHTOKEN"
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = "status=200"
$ws.Range("K4").Value = ""
$ws.Rows.Item(4).RowHeight = 30

# Mirror row 2's cell formatting (styles for D/E/F/G/H/I/J/K) onto rows 3 and 4
$ws.Range("D2:K2").Copy()
$ws.Range("D3:K3").PasteSpecial(-4122)
$ws.Range("D2:K2").Copy()
$ws.Range("D4:K4").PasteSpecial(-4122)
$ws.Range("B3:B4").WrapText = $true

# --- Column width adjustments ---
$ws.Columns.Item(2).ColumnWidth = 32.5
$ws.Columns.Item(4).ColumnWidth = 9.5
$ws.Columns.Item(6).ColumnWidth = 8.166666666666666
$ws.Columns.Item(8).ColumnWidth = 5.166666666666667
$ws.Columns.Item(10).ColumnWidth = 12.333333333333334

# --- Selection / view ---
$ws.Range("G4").Select()
